$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3717.1936
$ws.Range("J17").Value = 3717.1936
$ws.Range("L17").Value = 11151.5808
$ws.Range("N17").Value = -11487.5808
$ws.Range("H33").Value = 187.89473
$ws.Range("I33").Value = 227.64285
$ws.Range("J33").Value = 76.59999999999999
$ws.Range("K33").Value = 227.64285
$ws.Range("L33").Value = 76.59999999999999
$ws.Range("M33").Value = 1.35714999999999
$ws.Range("N33").Value = -534.6
$ws.Range("H138").Value = 5816.9697
$ws.Range("J138").Value = 7128.5625
$ws.Range("L138").Value = 21385.6875
$ws.Range("N138").Value = -31665.6875
$ws.Range("H139").Value = 99926.664
$ws.Range("J139").Value = 99926.664
$ws.Range("L139").Value = 99926.664
$ws.Range("N139").Value = -110206.664
$ws.Range("H140").Value = 68051
$ws.Range("J140").Value = 67834.44500000001
$ws.Range("L140").Value = 67834.44500000001
$ws.Range("N140").Value = -78194.44500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7005.2
$ws.Range("I32").Value = 4544.364
$ws.Range("K32").Value = 4544.364
$ws.Range("M32").Value = -4257.364
$ws.Range("H63").Value = 3100
$ws.Range("I63").Value = 3150
$ws.Range("K63").Value = 3150
$ws.Range("M63").Value = -2464
$ws.Range("H66").Value = 3100
$ws.Range("I66").Value = 3150
$ws.Range("K66").Value = 15750
$ws.Range("M66").Value = -12318
$ws.Range("H110").Value = 4402.5454
$ws.Range("I110").Value = 3824.5789
$ws.Range("K110").Value = 3824.5789
$ws.Range("M110").Value = -1779.5789
$ws.Range("H140").Value = 111933
$ws.Range("J140").Value = 111933
$ws.Range("L140").Value = 111933
$ws.Range("N140").Value = -122293

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H119").Value = 30000
$ws.Range("J119").Value = 30000
$ws.Range("L119").Value = 30000
$ws.Range("N119").Value = -39676
$ws.Range("H140").Value = 116627.19
$ws.Range("I140").Value = 51000
$ws.Range("J140").Value = 131771.92
$ws.Range("K140").Value = 51000
$ws.Range("L140").Value = 131771.92
$ws.Range("M140").Value = -45820
$ws.Range("N140").Value = -142131.92

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4260.8237
$ws.Range("I31").Value = 2930.617
$ws.Range("J31").Value = 7237.952
$ws.Range("K31").Value = 2930.617
$ws.Range("L31").Value = 7237.952
$ws.Range("M31").Value = -2635.617
$ws.Range("N31").Value = -7827.952
$ws.Range("H34").Value = 4260.8237
$ws.Range("I34").Value = 2930.617
$ws.Range("J34").Value = 7237.952
$ws.Range("K34").Value = 2930.617
$ws.Range("L34").Value = 7237.952
$ws.Range("M34").Value = -2728.617
$ws.Range("N34").Value = -7641.952
$ws.Range("H86").Value = 5657.6
$ws.Range("I86").Value = 5572
$ws.Range("K86").Value = 5572
$ws.Range("M86").Value = -4449
$ws.Range("H89").Value = 5657.6
$ws.Range("I89").Value = 5572
$ws.Range("K89").Value = 27860
$ws.Range("M89").Value = -22244
$ws.Range("H105").Value = 1620.3334
$ws.Range("I105").Value = 972.5
$ws.Range("J105").Value = 2138.6
$ws.Range("K105").Value = 972.5
$ws.Range("L105").Value = 2138.6
$ws.Range("M105").Value = 774.5
$ws.Range("N105").Value = -5632.6
$ws.Range("H109").Value = 58421.168
$ws.Range("J109").Value = 58421.168
$ws.Range("L109").Value = 58421.168
$ws.Range("N109").Value = -60501.168
$ws.Range("H141").Value = 80935.63
$ws.Range("I141").Value = 37999.668
$ws.Range("J141").Value = 88986.125
$ws.Range("K141").Value = 37999.668
$ws.Range("L141").Value = 88986.125
$ws.Range("M141").Value = -32819.668
$ws.Range("N141").Value = -99346.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 27070.857
$ws.Range("I74").Value = 26996
$ws.Range("J74").Value = 27083.334
$ws.Range("K74").Value = 80988
$ws.Range("L74").Value = 81250.00199999999
$ws.Range("M74").Value = -79927
$ws.Range("N74").Value = -83372.00199999999
$ws.Range("H77").Value = 27070.857
$ws.Range("I77").Value = 26996
$ws.Range("J77").Value = 27083.334
$ws.Range("K77").Value = 242964
$ws.Range("L77").Value = 243750.006
$ws.Range("M77").Value = -237660
$ws.Range("N77").Value = -254358.006
$ws.Range("H113").Value = 760
$ws.Range("J113").Value = 771.4286
$ws.Range("L113").Value = 2314.2858
$ws.Range("N113").Value = -6654.2858
$ws.Range("H131").Value = 2002.7162
$ws.Range("I131").Value = 2234.0908
$ws.Range("J131").Value = 1962.3175
$ws.Range("K131").Value = 6702.2724
$ws.Range("L131").Value = 5886.9525
$ws.Range("M131").Value = -1662.2724
$ws.Range("N131").Value = -15966.9525

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 52035.953
$ws.Range("I80").Value = 93296.17999999999
$ws.Range("J80").Value = 6649.7
$ws.Range("K80").Value = 93296.17999999999
$ws.Range("L80").Value = 6649.7
$ws.Range("M80").Value = -92298.17999999999
$ws.Range("N80").Value = -8645.700000000001
$ws.Range("H83").Value = 52035.953
$ws.Range("I83").Value = 93296.17999999999
$ws.Range("J83").Value = 6649.7
$ws.Range("K83").Value = 466480.9
$ws.Range("L83").Value = 33248.5
$ws.Range("M83").Value = -461488.9
$ws.Range("N83").Value = -43232.5
$ws.Range("H93").Value = 44545
$ws.Range("J93").Value = 44545
$ws.Range("L93").Value = 44545
$ws.Range("N93").Value = -48289
$ws.Range("H102").Value = 12965.941
$ws.Range("I102").Value = 14903.223
$ws.Range("K102").Value = 14903.223
$ws.Range("M102").Value = -13281.223
$ws.Range("H113").Value = 703778.1
$ws.Range("I113").Value = 1262880.6
$ws.Range("J113").Value = 4899.9165
$ws.Range("K113").Value = 1262880.6
$ws.Range("L113").Value = 4899.9165
$ws.Range("M113").Value = -1260710.6
$ws.Range("N113").Value = -9239.916499999999
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H126").Value = 4284.5713
$ws.Range("I126").Value = 2248.5
$ws.Range("K126").Value = 6745.5
$ws.Range("M126").Value = -4275.5
$ws.Range("H140").Value = 69780
$ws.Range("J140").Value = 69780
$ws.Range("L140").Value = 69780
$ws.Range("N140").Value = -80140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1646148.4
$ws.Range("I82").Value = 2605617.5
$ws.Range("J82").Value = 1344.1428
$ws.Range("K82").Value = 2605617.5
$ws.Range("L82").Value = 1344.1428
$ws.Range("M82").Value = -2605256.5
$ws.Range("N82").Value = -2066.1428
$ws.Range("H85").Value = 1646148.4
$ws.Range("I85").Value = 2605617.5
$ws.Range("J85").Value = 1344.1428
$ws.Range("K85").Value = 2605617.5
$ws.Range("L85").Value = 1344.1428
$ws.Range("M85").Value = -2604369.5
$ws.Range("N85").Value = -3840.1428

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 14000.667
$ws.Range("I62").Value = 11001
$ws.Range("J62").Value = 20000
$ws.Range("K62").Value = 11001
$ws.Range("L62").Value = 20000
$ws.Range("M62").Value = -10377
$ws.Range("N62").Value = -21248
$ws.Range("H65").Value = 14000.667
$ws.Range("I65").Value = 11001
$ws.Range("J65").Value = 20000
$ws.Range("K65").Value = 55005
$ws.Range("L65").Value = 100000
$ws.Range("M65").Value = -51885
$ws.Range("N65").Value = -106240
$ws.Range("H81").Value = 2089276.2
$ws.Range("I81").Value = 2609929
$ws.Range("J81").Value = 6665.5
$ws.Range("K81").Value = 5219858
$ws.Range("L81").Value = 13331
$ws.Range("M81").Value = -5218797
$ws.Range("N81").Value = -15453
$ws.Range("H84").Value = 2089276.2
$ws.Range("I84").Value = 2609929
$ws.Range("J84").Value = 6665.5
$ws.Range("K84").Value = 26099290
$ws.Range("L84").Value = 66655
$ws.Range("M84").Value = -26093986
$ws.Range("N84").Value = -77263
$ws.Range("H100").Value = 835348.8
$ws.Range("I100").Value = 1178165.5
$ws.Range("J100").Value = 2794
$ws.Range("K100").Value = 2356331
$ws.Range("L100").Value = 5588
$ws.Range("M100").Value = -2355790
$ws.Range("N100").Value = -6670
$ws.Range("H136").Value = 9405.662
$ws.Range("I136").Value = 4530.25
$ws.Range("J136").Value = 9970.928
$ws.Range("K136").Value = 13590.75
$ws.Range("L136").Value = 29912.784
$ws.Range("M136").Value = -11040.75
$ws.Range("N136").Value = -35012.784

